# Insert a new row of data at row 78 (pushes the existing rows 78..204 down
# to 79..205, extending the used range from A1:T204 to A1:T205), then
# populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(78).Insert()

$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 45219
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108002
$ws.Range("J78").Value = "Mango"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 250
$ws.Range("N78").Value = 10000
$ws.Range("O78").Value = 10000
$ws.Range("P78").Value = 10000
$ws.Range("Q78").Value = "`$/bandeja 4 kilos"
$ws.Range("R78").Value = "Brasil"
$ws.Range("S78").Value = 2500
$ws.Range("T78").Value = 4
